# Applies the "Updated cryptos list" data refresh (price / 1h-volume columns,
# plus the Kaspa/Filecoin row swap) described by the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", "69.512.75"),
    @("E2", "  +1.45%  "),
    @("D3", "3.888.96"),
    @("E3", "  +1.41%  "),
    @("E4", "  +0.04%  "),
    @("D5", "604.25"),
    @("E5", "  +0.64%  "),
    @("D6", "170.78"),
    @("E6", "  +4.32%  "),
    @("D7", "3.889.35"),
    @("E7", "  +1.46%  "),
    @("E8", "  -0.01%  "),
    @("E9", "  +1.07%  "),
    @("E10", "  +1.68%  "),
    @("E11", "  +1.07%  "),
    @("E12", "  +1.98%  "),
    @("E13", "  +4.76%  "),
    @("D14", "38.28"),
    @("E14", "  +3.92%  "),
    @("D15", "4.542.95"),
    @("E15", "  +1.40%  "),
    @("D16", "3.891.84"),
    @("E16", "  +1.62%  "),
    @("D17", "69.570.31"),
    @("E17", "  +1.25%  "),
    @("D18", "18.81"),
    @("E18", "  +9.73%  "),
    @("E19", "  +0.47%  "),
    @("E20", "  -0.82%  "),
    @("D21", "11.07"),
    @("E21", "  -0.84%  "),
    @("D22", "489.73"),
    @("E22", "  +0.66%  "),
    @("E23", "  +3.75%  "),
    @("D24", "0.0000166"),
    @("E24", "  +3.57%  "),
    @("D25", "85.36"),
    @("E25", "  +1.52%  "),
    @("D26", "2.30"),
    @("E26", "  +2.71%  "),
    @("D27", "12.39"),
    @("E27", "  +2.25%  "),
    @("D28", "10.11"),
    @("E28", "  +1.11%  "),
    @("E29", "  +0.23%  "),
    @("E30", "  +1.15%  "),
    @("D31", "4.038.94"),
    @("E31", "  +1.30%  "),
    @("E32", "  +1.42%  "),
    @("D33", "7.81"),
    @("E33", "  -0.55%  "),
    @("D34", "31.90"),
    @("E34", "  +0.23%  "),
    @("D35", "3.855.22"),
    @("E35", "  +1.95%  "),
    @("E36", "  -0.44%  "),
    @("D37", "3.42"),
    @("E37", "  +15.23%  "),
    @("B38", "Filecoin"),
    @("C38", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"),
    @("D38", "6.12"),
    @("E38", "  +4.00%  "),
    @("B39", "Kaspa"),
    @("C39", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"),
    @("D39", "0.143"),
    @("E39", "  +2.45%  "),
    @("E40", "  +0.57%  "),
    @("E41", "  +0.01%  "),
    @("E43", "  +4.47%  "),
    @("D44", "437.45"),
    @("E44", "  +1.99%  "),
    @("D45", "48.07"),
    @("E45", "  -0.85%  "),
    @("D46", "8.70"),
    @("E46", "  +3.33%  "),
    @("D48", "0.000278"),
    @("E48", "  +22.25%  "),
    @("E49", "  +2.41%  "),
    @("D50", "40.21"),
    @("E50", "  +3.82%  "),
    @("D51", "141.32"),
    @("E51", "  -0.92%  ")
)

foreach ($u in $updates) {
    $addr = $u[0]
    $val = $u[1]
    $rng = $ws.Range($addr)
    # Force text so Excel does not reinterpret numeric-looking strings
    # (e.g. "604.25") as real numbers - the source data is text throughout.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}
